$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 13.946
$ws.Range("C3").Value = 41.671999999999997

$ws.Range("B4").Value = 13.680999999999999

$ws.Range("D5").Value = 231.785
$ws.Range("E5").Value = 4492.3159999999998

$ws.Range("E6").Select()
